$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 33   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/19/2026  Through  1/25/2026"

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '0'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '***.*'
$ws.Range("N15").Value = -66.666666666666

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 8
$ws.Range("L16").Value = -60
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -90.697674418604

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -36.363636363636
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = -30
$ws.Range("L17").Value = 133.333333333333
$ws.Range("M17").Value = -30
$ws.Range("N17").Value = -65

# Row 18
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -66.666666666666
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = -60
$ws.Range("M18").Value = -71.428571428571
$ws.Range("N18").Value = -91.304347826087

# Row 19
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -7.5
$ws.Range("I19").Value = 33
$ws.Range("J19").Value = 37
$ws.Range("K19").Value = -10.810810810810
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -26.666666666666
$ws.Range("N19").Value = -31.25

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '***.*'
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 0
$ws.Range("L20").Value = -83.333333333333
$ws.Range("M20").Value = 0
$ws.Range("M20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N20").Value = -97.142857142857

# Row 21
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -21.428571428571
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = -20.289855072463
$ws.Range("I21").Value = 48
$ws.Range("J21").Value = 63
$ws.Range("K21").Value = -23.809523809523
$ws.Range("L21").Value = -22.580645161290
$ws.Range("M21").Value = -32.394366197183
$ws.Range("N21").Value = -72.093023255813

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = '0'
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = -66.666666666666
$ws.Range("M22").Value = 0
$ws.Range("M22").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = -20
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = 300

# Row 24
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = 1.785714285714
$ws.Range("I24").Value = 51
$ws.Range("J24").Value = 48
$ws.Range("K24").Value = 6.25
$ws.Range("L24").Value = 15.909090909090
$ws.Range("M24").Value = 6.25

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 17.857142857142
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 22
$ws.Range("K25").Value = 45.454545454545
$ws.Range("L25").Value = 68.421052631578

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 66.666666666666
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 16
$ws.Range("K26").Value = 62.5
$ws.Range("L26").Value = 8.333333333333

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = '0'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '***.*'

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 150
$ws.Range("L28").Value = 400

# Row 33
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = '0'

